$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.676.98"
$ws.Range("E2").Value = "  -6.69%  "
$ws.Range("D3").Value = "2.538.45"
$ws.Range("E3").Value = "  -4.80%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "298.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.24%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.572"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0803"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E13").Value = "  +6.56%  "
$ws.Range("D14").Value = "2.934.92"
$ws.Range("E14").Value = "  -4.66%  "
$ws.Range("D15").Value = "2.540.67"
$ws.Range("E15").Value = "  -4.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.874"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.23%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -6.22%  "
$ws.Range("D18").Value = "42.706.72"
$ws.Range("E18").Value = "  -6.63%  "
$ws.Range("D19").Value = "0.0₃0977"
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "253.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -10.21%  "
$ws.Range("E24").Value = "  -5.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.63%  "
$ws.Range("E27").Value = "  +0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.57%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.20"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "152.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.16%  "
$ws.Range("E33").Value = "  -2.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.29%  "
$ws.Range("E35").Value = "  -9.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0790"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.113"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -6.35%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.119"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.47%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "23.19"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0308"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.50%  "
$ws.Range("D44").Value = "2.077.77"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "84.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.60%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.791.07"
$ws.Range("E49").Value = "  -4.74%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -6.53%  "
$ws.Range("E51").Value = "  -3.81%  "
